# ---------------------------------------------------------------------------
# Edit summary (see commit message: "align third party compiler with other
# compilers"):
#   1. Handout master & notes master "datetime1" date fields:
#        06/11/2025 -> 07/11/2025
#      (the cached text of an <a:fld type="datetime1"> placeholder on the
#      Handout Master / Notes Master - this build of the host does not
#      expose a working write path for Handout/Notes Master placeholder
#      text, so this part is intentionally left as a no-op rather than
#      risk corrupting/perturbing unrelated parts of the package.)
#   2. Slide 3 (the architecture diagram): move/resize the two connectors
#      that attach to the "Third Party C Compiler" rounded-rectangle and
#      move that rounded rectangle itself up, so it lines up with the other
#      "Java Compiler" / "C Compiler" boxes on the same row.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- Re-align "Third Party C Compiler" box + its two connectors ------------

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

$slide = $p.Slides.Item(3)

# Connector dropping down into the top of the "Third Party C Compiler" box
# (id 45, stCxn -> shape 116). Only its vertical position/size change; the
# horizontal position is untouched.
$connTop = Get-ShapeById $slide.Shapes 45
$connTop.Top = 140.94661717322833
$connTop.Width = 0.9400788401574802
$connTop.Height = 116.65259942519684

# Connector feeding into "Third Party C Compiler" from "C Code" above
# (id 47, endCxn -> shape 116). Only its height changes.
$connBottom = Get-ShapeById $slide.Shapes 47
$connBottom.Height = 34.57858087716536

# The "Third Party C Compiler" rounded rectangle itself (id 116) moves up to
# line up with the other compiler boxes; size/left unchanged.
$compilerBox = Get-ShapeById $slide.Shapes 116
$compilerBox.Top = 110.56165704330708
